$wb = $excel.ActiveWorkbook

# Replace "Ready for handoff" with "In Translation" wherever it occurs
# (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4)
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NB: put the string literal on the left of -eq. PowerShell's -eq
        # coerces the right side to the left operand's type, and with a
        # boolean cell value on the right (e.g. the "True" cells) a
        # string-on-the-right comparison like `$cell.Value2 -eq "..."`
        # would coerce "Ready for handoff" to $true and false-match it.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the Status-related columns that previously fit "Ready for handoff".
# NOTE: the COM ColumnWidth setter here (like real Excel) quantizes the
# stored OOXML width to the nearest 1/6 after adding the fixed 5/6 padding
# term, so we dial in the character-width input that lands on the closest
# achievable grid point to the target 13.4101845877511.
$targetWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $targetWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $targetWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $targetWidth
